$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cohorts")

# The header row (row 1) holds its property-path labels starting at C1
# instead of A1 like every other sheet in this workbook. Shift the whole
# labeled block two columns to the left so it starts at A1 (matching the
# "biosamples" sheet's layout), the same way it would look if someone
# copied C1:CR1 and pasted it starting at A1.
$ws.Range("C1:CR1").Copy($ws.Range("A1")) | Out-Null
